$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 238 (shifts existing rows 238:327 down to 239:328)
$ws.Rows.Item(238).Insert()

# Populate the new row 238 with the new weekly price record
$ws.Cells.Item(238, 1).Value = 4
$ws.Cells.Item(238, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(238, 3).Value = 'Los Lagos'
$ws.Cells.Item(238, 4).Value = 44784
$ws.Cells.Item(238, 5).Value = 10
$ws.Cells.Item(238, 6).Value = 100112040
$ws.Cells.Item(238, 7).Value = 'Cilantro'
$ws.Cells.Item(238, 8).Value = 'Sin especificar'
$ws.Cells.Item(238, 9).Value = 'Primera'
$ws.Cells.Item(238, 10).Value = 70
$ws.Cells.Item(238, 11).Value = 12500
$ws.Cells.Item(238, 12).Value = 12500
$ws.Cells.Item(238, 13).Value = 12500
$ws.Cells.Item(238, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(238, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(238, 16).Value = 347
$ws.Cells.Item(238, 17).Value = 36
$ws.Cells.Item(238, 18).Value = 'Hortaliza'
